# Fill in the previously-missing "Verdi" (F column) values with 10 for the
# rows where that figure had not yet been entered, then leave the
# selection where the author left it when they saved (cell J73).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(30, 31, 46, 47, 48, 49, 56, 57, 58, 59)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 10
}

$ws.Range("J73").Select()
